# Scheduled market-data refresh: update price/profit columns (H:N) on each profession sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 137
$ws.Range("H137").Value = 2828.16
$ws.Range("I137").Value = 2040.55
$ws.Range("J137").Value = 5978.6
$ws.Range("K137").Value = 6121.65
$ws.Range("L137").Value = 17935.8
$ws.Range("M137").Value = -3571.65
$ws.Range("N137").Value = -23035.8
# Row 138
$ws.Range("H138").Value = 3909.2104
$ws.Range("I138").Value = 9282
$ws.Range("J138").Value = 2901.8125
$ws.Range("K138").Value = 27846
$ws.Range("L138").Value = 8705.4375
$ws.Range("M138").Value = -22706
$ws.Range("N138").Value = -18985.4375

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2127.5417
$ws.Range("I61").Value = 2209.9412
$ws.Range("J61").Value = 1927.4286
$ws.Range("K61").Value = 2209.9412
$ws.Range("L61").Value = 1927.4286
$ws.Range("M61").Value = -1997.9412
$ws.Range("N61").Value = -2351.4286
# Row 74
$ws.Range("H74").Value = 2786.1
$ws.Range("I74").Value = 1199.9333
$ws.Range("K74").Value = 1199.9333
$ws.Range("M74").Value = -325.9332999999999
# Row 77
$ws.Range("H77").Value = 2786.1
$ws.Range("I77").Value = 1199.9333
$ws.Range("K77").Value = 5999.666499999999
$ws.Range("M77").Value = -1631.666499999999
# Row 136
$ws.Range("H136").Value = 2127.5417
$ws.Range("I136").Value = 2209.9412
$ws.Range("J136").Value = 1927.4286
$ws.Range("K136").Value = 6629.823600000001
$ws.Range("L136").Value = 5782.2858
$ws.Range("M136").Value = -4079.823600000001
$ws.Range("N136").Value = -10882.2858

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 665.625
$ws.Range("I80").Value = 873.2
$ws.Range("J80").Value = 319.66666
$ws.Range("K80").Value = 873.2
$ws.Range("L80").Value = 319.66666
$ws.Range("M80").Value = 124.8
$ws.Range("N80").Value = -2315.66666
# Row 83
$ws.Range("H83").Value = 665.625
$ws.Range("I83").Value = 873.2
$ws.Range("J83").Value = 319.66666
$ws.Range("K83").Value = 4366
$ws.Range("L83").Value = 1598.3333
$ws.Range("M83").Value = 626
$ws.Range("N83").Value = -11582.3333
# Row 99
$ws.Range("H99").Value = 3131.95
$ws.Range("I99").Value = 3581
$ws.Range("J99").Value = 2583.111
$ws.Range("K99").Value = 3581
$ws.Range("L99").Value = 2583.111
$ws.Range("M99").Value = -2083
$ws.Range("N99").Value = -5579.111
# Row 134
$ws.Range("H134").Value = 2857.0476
$ws.Range("I134").Value = 2535.4
$ws.Range("J134").Value = 3661.1667
$ws.Range("K134").Value = 7606.200000000001
$ws.Range("L134").Value = 10983.5001
$ws.Range("M134").Value = -5071.200000000001
$ws.Range("N134").Value = -16053.5001

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 492.2
$ws.Range("I16").Value = 492.2
$ws.Range("K16").Value = 492.2
$ws.Range("M16").Value = -205.2
# Row 31
$ws.Range("H31").Value = 3830.0688
$ws.Range("I31").Value = 2658
$ws.Range("J31").Value = 5748
$ws.Range("K31").Value = 2658
$ws.Range("L31").Value = 5748
$ws.Range("M31").Value = -2363
$ws.Range("N31").Value = -6338
# Row 34
$ws.Range("H34").Value = 3830.0688
$ws.Range("I34").Value = 2658
$ws.Range("J34").Value = 5748
$ws.Range("K34").Value = 2658
$ws.Range("L34").Value = 5748
$ws.Range("M34").Value = -2456
$ws.Range("N34").Value = -6152
# Row 58
$ws.Range("H58").Value = 3211.4167
$ws.Range("I58").Value = 1188.7273
$ws.Range("J58").Value = 4922.923
$ws.Range("K58").Value = 1188.7273
$ws.Range("L58").Value = 4922.923
$ws.Range("M58").Value = -985.7273
$ws.Range("N58").Value = -5328.923
# Row 113
$ws.Range("H113").Value = 492.2
$ws.Range("I113").Value = 492.2
$ws.Range("K113").Value = 492.2
$ws.Range("M113").Value = 1677.8
# Row 122
$ws.Range("H122").Value = 2550.158
$ws.Range("I122").Value = 2550.158
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7650.474
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5200.474
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 3761.55
$ws.Range("I132").Value = 2826.125
$ws.Range("J132").Value = 7503.25
$ws.Range("K132").Value = 8478.375
$ws.Range("L132").Value = 22509.75
$ws.Range("M132").Value = -5948.375
$ws.Range("N132").Value = -27569.75
# Row 134
$ws.Range("H134").Value = 2211.15
$ws.Range("I134").Value = 1113.5834
$ws.Range("J134").Value = 3857.5
$ws.Range("K134").Value = 3340.7502
$ws.Range("L134").Value = 11572.5
$ws.Range("M134").Value = -805.7501999999999
$ws.Range("N134").Value = -16642.5
# Row 136
$ws.Range("H136").Value = 3211.4167
$ws.Range("I136").Value = 1188.7273
$ws.Range("J136").Value = 4922.923
$ws.Range("K136").Value = 3566.1819
$ws.Range("L136").Value = 14768.769
$ws.Range("M136").Value = -1016.1819
$ws.Range("N136").Value = -19868.769

$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 3334000
$ws.Range("J46").Value = 3334000
$ws.Range("L46").Value = 10002000
$ws.Range("N46").Value = -10002182
# Row 68
$ws.Range("H68").Value = 1015.5714
$ws.Range("I68").Value = 1062
$ws.Range("J68").Value = 997
$ws.Range("K68").Value = 3186
$ws.Range("L68").Value = 2991
$ws.Range("M68").Value = -2375
$ws.Range("N68").Value = -4613
# Row 71
$ws.Range("H71").Value = 1015.5714
$ws.Range("I71").Value = 1062
$ws.Range("J71").Value = 997
$ws.Range("K71").Value = 9558
$ws.Range("L71").Value = 8973
$ws.Range("M71").Value = -5502
$ws.Range("N71").Value = -17085

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 1946.2413
$ws.Range("I132").Value = 891.1579
$ws.Range("J132").Value = 3950.9
$ws.Range("K132").Value = 2673.4737
$ws.Range("L132").Value = 11852.7
$ws.Range("M132").Value = -143.4737
$ws.Range("N132").Value = -16912.7

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 500
$ws.Range("I7").Value = 500
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 500
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -388
$ws.Range("N7").ClearContents()
# Row 93
$ws.Range("H93").Value = 2649.889
$ws.Range("I93").Value = 2649.889
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2649.889
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -1401.889
$ws.Range("N93").ClearContents()
# Row 122
$ws.Range("H122").Value = 4834.1904
$ws.Range("I122").Value = 4229.7
$ws.Range("J122").Value = 5383.727
$ws.Range("K122").Value = 12689.1
$ws.Range("L122").Value = 16151.181
$ws.Range("M122").Value = -10239.1
$ws.Range("N122").Value = -21051.181
# Row 126
$ws.Range("H126").Value = 500
$ws.Range("I126").Value = 500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 1500
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 970
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 5135.2856
$ws.Range("I132").Value = 4203.8184
$ws.Range("K132").Value = 12611.4552
$ws.Range("M132").Value = -10081.4552
# Row 136
$ws.Range("H136").Value = 2853.4666
$ws.Range("J136").Value = 3132.1667
$ws.Range("L136").Value = 9396.500100000001
$ws.Range("N136").Value = -14496.5001

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2206.16
$ws.Range("I122").Value = 2286.318
$ws.Range("J122").Value = 1618.3334
$ws.Range("K122").Value = 6858.954000000001
$ws.Range("L122").Value = 4855.0002
$ws.Range("M122").Value = -4408.954000000001
$ws.Range("N122").Value = -9755.0002
# Row 132
$ws.Range("H132").Value = 740.4167
$ws.Range("I132").Value = 731.78125
$ws.Range("J132").Value = 809.5
$ws.Range("K132").Value = 2195.34375
$ws.Range("L132").Value = 2428.5
$ws.Range("M132").Value = 334.65625
$ws.Range("N132").Value = -7488.5
# Row 136
$ws.Range("H136").Value = 3636.5
$ws.Range("I136").Value = 1885.1724
$ws.Range("J136").Value = 9279.666999999999
$ws.Range("K136").Value = 5655.5172
$ws.Range("L136").Value = 27839.001
$ws.Range("M136").Value = -3105.5172
$ws.Range("N136").Value = -32939.001
